# Update computed profit/price figures on several sheets to reflect
# refreshed market data (scheduled runner sync).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 775.2
$ws.Range("I6").Value = 231.5
$ws.Range("J6").Value = 2950
$ws.Range("K6").Value = 694.5
$ws.Range("L6").Value = 8850
$ws.Range("M6").Value = -582.5
$ws.Range("N6").Value = -9074
$ws.Range("H31").Value = 9515.071
$ws.Range("I31").Value = 9515.071
$ws.Range("K31").Value = 28545.213
$ws.Range("M31").Value = -28315.213
$ws.Range("H33").Value = 150.12903
$ws.Range("I33").Value = 159.96297
$ws.Range("J33").Value = 83.75
$ws.Range("K33").Value = 159.96297
$ws.Range("L33").Value = 83.75
$ws.Range("M33").Value = 69.03702999999999
$ws.Range("N33").Value = -541.75
$ws.Range("H93").Value = 255000
$ws.Range("I93").Value = 10000
$ws.Range("J93").Value = 500000
$ws.Range("K93").Value = 10000
$ws.Range("L93").Value = 500000
$ws.Range("M93").Value = -7504
$ws.Range("N93").Value = -504992
$ws.Range("H98").Value = 4086.4243
$ws.Range("I98").Value = 2910.7407
$ws.Range("J98").Value = 9377
$ws.Range("K98").Value = 2910.7407
$ws.Range("L98").Value = 9377
$ws.Range("M98").Value = -1412.7407
$ws.Range("N98").Value = -12373
$ws.Range("H122").Value = 4086.4243
$ws.Range("I122").Value = 2910.7407
$ws.Range("J122").Value = 9377
$ws.Range("K122").Value = 8732.222099999999
$ws.Range("L122").Value = 28131
$ws.Range("M122").Value = -6282.222099999999
$ws.Range("N122").Value = -33031
$ws.Range("H129").Value = 856.58
$ws.Range("J129").Value = 933.13794
$ws.Range("L129").Value = 2799.41382
$ws.Range("N129").Value = -12799.41382
$ws.Range("H132").Value = 2041.5834
$ws.Range("I132").Value = 1713.1072
$ws.Range("K132").Value = 5139.321599999999
$ws.Range("M132").Value = -2609.321599999999
$ws.Range("H134").Value = 77333.336
$ws.Range("I134").Value = 75000
$ws.Range("J134").Value = 77800
$ws.Range("K134").Value = 75000
$ws.Range("L134").Value = 77800
$ws.Range("M134").Value = -69930
$ws.Range("N134").Value = -87940
$ws.Range("H137").Value = 3248.766
$ws.Range("I137").Value = 1876.4073
$ws.Range("J137").Value = 3989.84
$ws.Range("K137").Value = 5629.2219
$ws.Range("L137").Value = 11969.52
$ws.Range("M137").Value = -3079.2219
$ws.Range("N137").Value = -17069.52
$ws.Range("H138").Value = 4043.5293
$ws.Range("I138").Value = 1421.4783
$ws.Range("J138").Value = 5383.689
$ws.Range("K138").Value = 4264.4349
$ws.Range("L138").Value = 16151.067
$ws.Range("M138").Value = 875.5650999999998
$ws.Range("N138").Value = -26431.067

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 25389.545
$ws.Range("J123").Value = 25389.545
$ws.Range("L123").Value = 25389.545
$ws.Range("N123").Value = -35189.545

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 30000
$ws.Range("J17").Value = 30000
$ws.Range("L17").Value = 30000
$ws.Range("N17").Value = -30348
$ws.Range("H31").Value = 4281.814
$ws.Range("I31").Value = 1915.3846
$ws.Range("J31").Value = 5307.2666
$ws.Range("K31").Value = 1915.3846
$ws.Range("L31").Value = 5307.2666
$ws.Range("M31").Value = -1620.3846
$ws.Range("N31").Value = -5897.2666
$ws.Range("H34").Value = 4281.814
$ws.Range("I34").Value = 1915.3846
$ws.Range("J34").Value = 5307.2666
$ws.Range("K34").Value = 1915.3846
$ws.Range("L34").Value = 5307.2666
$ws.Range("M34").Value = -1713.3846
$ws.Range("N34").Value = -5711.2666
$ws.Range("H135").Value = 71930.836
$ws.Range("J135").Value = 71930.836
$ws.Range("L135").Value = 71930.836
$ws.Range("N135").Value = -82070.836

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 5000
$ws.Range("J49").Value = 5000
$ws.Range("L49").Value = 15000
$ws.Range("N49").Value = -15312
$ws.Range("H108").Value = 3283.923
$ws.Range("I108").Value = 797
$ws.Range("K108").Value = 2391
$ws.Range("M108").Value = 489
$ws.Range("H118").Value = 3311.2727
$ws.Range("I118").Value = 2000
$ws.Range("J118").Value = 3442.4
$ws.Range("K118").Value = 6000
$ws.Range("L118").Value = 10327.2
$ws.Range("M118").Value = -4757
$ws.Range("N118").Value = -12813.2
$ws.Range("H129").Value = 3571927.8
$ws.Range("J129").Value = 16668000
$ws.Range("L129").Value = 50004000
$ws.Range("N129").Value = -50014000

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3071.875
$ws.Range("I40").Value = 2736.25
$ws.Range("K40").Value = 2736.25
$ws.Range("M40").Value = -2600.25
$ws.Range("H46").Value = 1785.5714
$ws.Range("I46").Value = 1500
$ws.Range("J46").Value = 1999.75
$ws.Range("K46").Value = 1500
$ws.Range("L46").Value = 1999.75
$ws.Range("M46").Value = -1312
$ws.Range("N46").Value = -2375.75
$ws.Range("H122").Value = 8699638
$ws.Range("I122").Value = 2288.5
$ws.Range("J122").Value = 13338224
$ws.Range("K122").Value = 6865.5
$ws.Range("L122").Value = 40014672
$ws.Range("M122").Value = -4415.5
$ws.Range("N122").Value = -40019572

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 20128
$ws.Range("J75").Value = 20128
$ws.Range("L75").Value = 20128
$ws.Range("N75").Value = -22000
$ws.Range("H78").Value = 20128
$ws.Range("J78").Value = 20128
$ws.Range("L78").Value = 60384
$ws.Range("N78").Value = -69744
$ws.Range("H81").Value = 41333.17
$ws.Range("I81").Value = 57893.1
$ws.Range("J81").Value = 4533.3335
$ws.Range("K81").Value = 115786.2
$ws.Range("L81").Value = 9066.666999999999
$ws.Range("M81").Value = -114725.2
$ws.Range("N81").Value = -11188.667
$ws.Range("H84").Value = 41333.17
$ws.Range("I84").Value = 57893.1
$ws.Range("J84").Value = 4533.3335
$ws.Range("K84").Value = 578931
$ws.Range("L84").Value = 45333.335
$ws.Range("M84").Value = -573627
$ws.Range("N84").Value = -55941.335
$ws.Range("H123").Value = 23803.166
$ws.Range("J123").Value = 23803.166
$ws.Range("L123").Value = 23803.166
$ws.Range("N123").Value = -33603.166
$ws.Range("H125").Value = 48672.25
$ws.Range("J125").Value = 48672.25
$ws.Range("L125").Value = 48672.25
$ws.Range("N125").Value = -58512.25
$ws.Range("H128").Value = 51715
$ws.Range("J128").Value = 51715
$ws.Range("L128").Value = 51715
$ws.Range("N128").Value = -61675
$ws.Range("H138").Value = 40984.8
$ws.Range("J138").Value = 40984.8
$ws.Range("L138").Value = 40984.8
$ws.Range("N138").Value = -51264.8
